# Updates cryptos list (prices / 1h volume%) to match the Dec 2, 2023
# GitHub Actions refresh. Price cells that look like plain numbers are
# written with a leading apostrophe so Excel keeps them as text (matching
# the original inlineStr cells) instead of silently converting to numeric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.772.91'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.100.84'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''226.55'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").Value = '''62.17'
$ws.Range("E7").Value = '  +3.12%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  +6.40%  '
$ws.Range("D13").Value = '2.412.31'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '''21.90'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '''0.801'
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '2.131.48'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("D18").Value = '38.805.40'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").Value = '''71.59'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '''6.04'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").Value = '''227.28'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '''2.55'
$ws.Range("E24").Value = '  +5.05%  '
$ws.Range("D25").Value = '''2.30'
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("D26").Value = '''9.64'
$ws.Range("E26").Value = '  +2.23%  '
$ws.Range("D27").Value = '''170.66'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("D30").Value = '''19.33'
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("E31").Value = '  +9.15%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''4.77'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("B35").Value = 'THORChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D35").Value = '''7.13'
$ws.Range("E35").Value = '  +12.97%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").Value = '''2.36'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").Value = '''17.97'
$ws.Range("E41").Value = '  +3.41%  '
$ws.Range("D42").Value = '''101.15'
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D43").Value = '1.525.75'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '''1.20'
$ws.Range("E44").Value = '  +7.59%  '
$ws.Range("D45").Value = '''2.81'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''7.79'
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.0915'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("E48").Value = '  +5.00%  '
$ws.Range("D49").Value = '''4.15'
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").Value = '2.299.86'
$ws.Range("E51").Value = '  +1.09%  '
